$d = $word.ActiveDocument

# Insert the new content (heading, intro line, and the flow-chart list) as a
# block of brand-new paragraphs before everything that is currently in the
# document body.
$insertionPoint = $d.Range(0, 0)
$newContentXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
            </w:pPr>
            <w:r>
              <w:t>Monday, May 18, 15</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>Flow chart for path finder.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>Find median = all &lt; 400?</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>Yes, FORWARD 10ms, GOTO 1</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>Median spread is balanced and &gt;= 3 sensors wide?</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>Yes, cross found, stop.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>Median to right?</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>Yes, RIGHT 10ms, GOTO 2</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>Median to left?</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>Yes, Left 10ms, GOTO 1</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>All &lt; 400?</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
            </w:pPr>
            <w:r>
              <w:t>Yes, end, stop.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$insertionPoint.InsertXML($newContentXml)

# Paragraphs 3..12 are the flow-chart steps; turn them into a genuine
# multilevel numbered list (1), 2), ... / a., b., ... ) sharing one list
# definition, alternating between outline level 0 and level 1.
$first = $d.Paragraphs(3).Range
$first.ListFormat.ApplyNumberDefault()
$listTemplate = $first.ListFormat.ListTemplate

$level1 = $listTemplate.ListLevels.Item(1)
$level1.NumberFormat = "%1)"

$level2 = $listTemplate.ListLevels.Item(2)
$level2.NumberStyle = 4
$level2.NumberFormat = "%2."

$level3 = $listTemplate.ListLevels.Item(3)
$level3.NumberStyle = 2
$level3.NumberFormat = "%3."

$level4 = $listTemplate.ListLevels.Item(4)
$level4.NumberFormat = "%4."

$level5 = $listTemplate.ListLevels.Item(5)
$level5.NumberStyle = 4
$level5.NumberFormat = "%5."

$level6 = $listTemplate.ListLevels.Item(6)
$level6.NumberStyle = 2
$level6.NumberFormat = "%6."

$level7 = $listTemplate.ListLevels.Item(7)
$level7.NumberFormat = "%7."

$level8 = $listTemplate.ListLevels.Item(8)
$level8.NumberStyle = 4
$level8.NumberFormat = "%8."

$level9 = $listTemplate.ListLevels.Item(9)
$level9.NumberStyle = 2
$level9.NumberFormat = "%9."

$levelForParagraph = @(1, 2, 1, 2, 1, 2, 1, 2, 1, 2)
for ($i = 3; $i -le 12; $i++) {
    $para = $d.Paragraphs($i).Range
    if ($i -gt 3) {
        $para.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true, 0, $false, 2)
    }
    $para.ListFormat.ListLevelNumber = $levelForParagraph[$i - 3]
}

Write-Host "Inserted new flow-chart section with $($d.Paragraphs.Count) total paragraphs."
